$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'The candidate has a strong technical background in frontend development, as evidenced by the projects ''LLMGuard'', ''Literature Society IITJ Website'', and skills in ReactJS, JavaScript, and Django. The applicant has also demonstrated proficiency in computer vision and image processing, as shown by the ''CloudPhysician''s Vital Extraction Challenge'' project. However, there are some skills required by the company, such as NextJS and CSS, that the applicant does not possess. Despite this, the candidate''s experience and skills showcase an ability to work in a team and ensure consistent web design and user experience. The applicant''s ability to create user-friendly interfaces and backend systems aligns well with the job description, making the candidate a suitable choice for the position.'
$ws.Range("E2").Value = 'The applicant has demonstrated strong communication skills, willingness to learn from Japanese work culture, and adaptability to different cultures. The applicant''s proactive approach to skill improvement and seeking feedback, as well as the desire to work in Japan to align with shared values of hard work and dedication, exhibits a strong personality fit for the company''s requirements. The applicant also expresses a clear career plan and expectations from the company, demonstrating a high level of motivation and commitment. Overall, the applicant''s responses reflect a positive and adaptable personality suitable for working in Japan.'
$ws.Range("C3").Value = 'The applicant has a strong suit in frontend and backend technologies including ReactJS, NodeJS, and CSS, which are directly in line with the company''s requirements. Additionally, the applicant''s projects demonstrate proficiency in using various technologies such as Flutter, Dart, Firebase, and WebRTC, showcasing versatility and an ability to adapt to different tech stacks. However, there is a lack of experience or projects related to NextJS, which is a requirement by the company. Despite this, the applicant''s wide array of skills and demonstrated project work reflect a high level of technical competence in web development, making them well-suited for the job.'
$ws.Range("E3").Value = 'The applicant has displayed strong communication skills, adaptability, and a collaborative approach to problem-solving, as evidenced by their responses. Their willingness to work in Japan is evident from their keen interest in the country''s work culture, their openness to learning a new language, and their willingness to adapt to a new culture. They also express a desire to work on cutting-edge AI/ML technologies, aligning with the company''s focus. Overall, the applicant''s responses indicate a good fit for the company''s soft skill requirements and a strong willingness to work in Japan.'
$ws.Range("C4").Value = 'The applicant''s projects demonstrate a strong technical background, particularly in Python, PyTorch, and ReactJS, which aligns with the company''s requirements for proficiency in ReactJS and frontend development. The ''Multi Model Data Analysis for Annotation of Human Activities'' project showcases the applicant''s skills in Python, PyTorch, and TensorFlow, which are beneficial for AI and machine learning roles. However, the applicant''s lack of experience with JavaScript, CSS, and NextJS could be seen as a gap in the required skills for the frontend engineer intern role. Despite this, the candidate''s proficiency in other relevant technologies and their ability to contribute to cutting-edge AI/ML technologies make them a strong contender for the position.'
$ws.Range("E4").Value = 'The applicant''s responses indicate a willingness to adapt to new cultures and work in Japan, as well as an interest in being part of the work culture and learning a new language. The applicant''s emphasis on teamwork, problem-solving, and adaptability align with the soft skills required by the company. Furthermore, their openness to working in Japan and their desire to learn new skills demonstrate a positive attitude and enthusiasm for the opportunity. Overall, the applicant''s responses reflect a strong potential for success in a new work environment, earning a score of 4 out of 5.'
$ws.Range("C5").Value = 'The applicant has demonstrated strong skills in frontend development, particularly in JavaScript, HTML, and CSS through their projects such as the Video Conferencing Project and Alcheringa Pass Portal. They have also shown the ability to work on a team project, which aligns with the requirement for excellent communication skills and teamwork in the job description. However, the applicant lacks experience with ReactJS and NextJS, which are required skills for the position. Therefore, although the applicant possesses relevant skills in frontend development and has demonstrated teamwork, the lack of experience with ReactJS and NextJS is a limiting factor in their suitability for the job.'
$ws.Range("E5").Value = 'The applicant has expressed a strong interest in being a part of the work culture in Japan and has shown openness to adapting to a new culture. They have also emphasized their ability to work well in a team, which aligns with the company''s requirement for teamwork. Additionally, the applicant has highlighted their interest in learning new skills and contributing to cutting-edge technologies, indicating a proactive and adaptable personality. However, the applicant''s preference to not work alone and the need to improve in acquiring new skills may pose some challenges in a new work environment. Overall, the applicant''s openness to adapt to a new culture, willingness to learn, and interest in collaborative work make them a suitable candidate for the role in Japan.'
